# Update simulation parameter / result values in row 5 of Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K5").Value = 0.00001
$ws.Range("L5").Value = 0.001

$ws.Range("Q5").Value = 126
$ws.Range("R5").Value = 83

$ws.Range("T5").Value = 0.1178742211419419
$ws.Range("U5").Value = 0.118437665648738

$ws.Range("W5").Value = 0.8619405466887732
$ws.Range("X5").Value = 0.8655080918800576
